# Insert two new data rows into the historic-change table to include the
# "Ragwitz et al. (2023)" / Germany "Nachfrage+Tech" scenario entries
# (one for the Energy sector, one for the Industry sector).
#
# Before the edit the sheet has 97 data rows (A1:K97). The new rows are
# inserted so that:
#   - a new Energy row lands at row 75 (pushing the former rows 75-85 down by one)
#   - a new Industry row lands at row 86 (pushing the former rows 86-97 down by one more)
# resulting in a sheet with 99 rows (A1:K99).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new "Energy" row at row 75 -------------------------------
$ws.Rows.Item(75).Insert()

$ws.Cells.Item(75, 1).Value = "Ragwitz et al. (2023)"
$ws.Cells.Item(75, 2).Value = "DE"
$ws.Cells.Item(75, 3).Value = "Nachfrage+Tech"
$ws.Cells.Item(75, 4).Value = "Energy"
$ws.Cells.Item(75, 5).Value = "final energy demand per capita and year"
$ws.Cells.Item(75, 6).Value = 2045
$ws.Cells.Item(75, 7).Value = "GJ/cap/year"
$ws.Cells.Item(75, 8).Value = 60.08
$ws.Cells.Item(75, 9).Value = 102.96
$ws.Cells.Item(75, 10).Value = -41.64724164724164
$ws.Cells.Item(75, 11).Value = 439

# --- Insert the new "Industry" row at row 86 (post first insert) --------
$ws.Rows.Item(86).Insert()

$ws.Cells.Item(86, 1).Value = "Ragwitz et al. (2023)"
$ws.Cells.Item(86, 2).Value = "DE"
$ws.Cells.Item(86, 3).Value = "Nachfrage+Tech"
$ws.Cells.Item(86, 4).Value = "Industry"
$ws.Cells.Item(86, 5).Value = "final energy demand per capita and year | industry"
$ws.Cells.Item(86, 6).Value = 2045
$ws.Cells.Item(86, 7).Value = "GJ/cap/year"
$ws.Cells.Item(86, 8).Value = 24.47
$ws.Cells.Item(86, 9).Value = 29.12
$ws.Cells.Item(86, 10).Value = -15.9684065934066
$ws.Cells.Item(86, 11).Value = 438
